$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to be treated as text so numeric-looking strings
    # (e.g. "215.17", "0.0880") are not coerced into actual numbers,
    # then drop back to the default "Normal" style so no stray style
    # index is left behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "27.608.58"
$ws.Range("E2").Value = "  -1.49%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "1.668.63"
$ws.Range("E3").Value = "  -3.06%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.20%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "215.17"
$ws.Range("E5").Value = "  -1.59%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -2.29%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.26%  "

# Row 8 - Solana
Set-TextValue $ws.Range("D8") "23.74"
$ws.Range("E8").Value = "  -1.37%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -0.62%  "

# Row 10 - Dogecoin
Set-TextValue $ws.Range("D10") "0.0623"
$ws.Range("E10").Value = "  -1.43%  "

# Row 11 - TRON
Set-TextValue $ws.Range("D11") "0.0880"
$ws.Range("E11").Value = "  -2.02%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D12") "1.903.35"
$ws.Range("E12").Value = "  -3.21%  "

# Row 13 - WrappedEther
Set-TextValue $ws.Range("D13") "1.673.50"
$ws.Range("E13").Value = "  -2.85%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -2.94%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +0.57%  "

# Row 16 - Litecoin
Set-TextValue $ws.Range("D16") "66.30"
$ws.Range("E16").Value = "  -1.60%  "

# Row 17 - was WrappedBTC, now BitcoinCash (rows 17/18 swapped content)
$ws.Range("B17").Value = "BitcoinCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D17") "243.99"
$ws.Range("E17").Value = "  +0.74%  "

# Row 18 - was BitcoinCash, now WrappedBTC
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D18") "27.593.59"
$ws.Range("E18").Value = "  -1.44%  "

# Row 19 - ShibaInu
$ws.Range("E19").Value = "  -3.34%  "

# Row 20 - Chainlink
Set-TextValue $ws.Range("D20") "7.63"
$ws.Range("E20").Value = "  -3.41%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.06%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -2.84%  "

# Row 23 - Avalanche
Set-TextValue $ws.Range("D23") "9.30"
$ws.Range("E23").Value = "  -3.78%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -4.41%  "

# Row 25 - Monero
Set-TextValue $ws.Range("D25") "146.78"
$ws.Range("E25").Value = "  -1.22%  "

# Row 26 - Cosmos
Set-TextValue $ws.Range("D26") "7.21"
$ws.Range("E26").Value = "  -3.66%  "

# Row 27 - EthereumClassic
Set-TextValue $ws.Range("D27") "16.46"
$ws.Range("E27").Value = "  -1.26%  "

# Row 28 - BinanceUSD
Set-TextValue $ws.Range("D28") "0.999"
$ws.Range("E28").Value = "  -0.30%  "

# Row 29 - Stellar
$ws.Range("E29").Value = "  -2.22%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +2.73%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  -1.50%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -2.37%  "

# Row 33 - Maker
Set-TextValue $ws.Range("D33") "1.467.28"
$ws.Range("E33").Value = "  -1.77%  "

# Row 34 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D34") "3.12"
$ws.Range("E34").Value = "  -4.51%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -4.79%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  -1.61%  "

# Row 37 - ARBITRUM
Set-TextValue $ws.Range("D37") "0.932"
$ws.Range("E37").Value = "  -2.06%  "

# Row 38 - was ImmutableX, now VeChain (rows 38/39 swapped content)
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D38") "0.0172"
$ws.Range("E38").Value = "  -0.87%  "

# Row 39 - was VeChain, now ImmutableX
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D39") "0.575"
$ws.Range("E39").Value = "  -4.88%  "

# Row 40 - Aave
Set-TextValue $ws.Range("D40") "69.58"
$ws.Range("E40").Value = "  -1.22%  "

# Row 41 - WEMIXToken
$ws.Range("E41").Value = "  -5.16%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  -0.19%  "

# Row 43 - FraxShare
Set-TextValue $ws.Range("D43") "5.42"
$ws.Range("E43").Value = "  -6.79%  "

# Row 44 - MXToken
$ws.Range("E44").Value = "  -3.65%  "

# Row 45 - RocketPoolETH
Set-TextValue $ws.Range("D45") "1.811.82"
$ws.Range("E45").Value = "  -3.19%  "

# Row 46 - TrustWalletToken
Set-TextValue $ws.Range("D46") "0.789"
$ws.Range("E46").Value = "  -0.72%  "

# Row 47 - RenderToken
$ws.Range("E47").Value = "  -2.56%  "

# Row 48 - Quant
Set-TextValue $ws.Range("D48") "89.44"
$ws.Range("E48").Value = "  -1.47%  "

# Row 49 - BabyDogeCoin
Set-TextValue $ws.Range("D49") "0.0₆0108"
$ws.Range("E49").Value = "  -3.61%  "

# Row 50 - Algorand
$ws.Range("E50").Value = "  -1.79%  "

# Row 51 - EnergySwap
Set-TextValue $ws.Range("D51") "7.87"
$ws.Range("E51").Value = "  -4.12%  "
